$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header "stcode11" in B1 to "code11"
$ws.Range("B1").Value = "code11"

# Reset the selection to B1 as seen in the updated workbook
$ws.Range("B1").Select()
